{"js": "// Build the \"CONTRATO DE MUTUO CONVERTIBLE\" template body.\n// Each entry is [styleBuiltIn-or-null, text]; null keeps the Normal style.\nconst items = [\n  [\"heading1\", \"CONTRATO DE MUTUO CONVERTIBLE\"],\n  [null, \"Entre ______________________ (el \\u201CPrestamista\\u201D) y la Compa\\u00F1\\u00EDa en formaci\\u00F3n (el \\u201CMutuario\\u201D), con fecha ___ de __________ de 2025.\"],\n  [\"heading2\", \"1. Monto y Desembolso\"],\n  [null, \"El Prestamista entrega USD __________ al Mutuario.\"],\n  [\"heading2\", \"2. Conversi\\u00F3n Autom\\u00E1tica\"],\n  [null, \"Al cerrar una ronda \\u2265 USD ________, el mutuo se convierte en acciones preferidas con 20 % de descuento o valuation cap USD ________, lo que convenga al Prestamista.\"],\n  [\"heading2\", \"3. Inter\\u00E9s y Ajuste\"],\n  [null, \"Inter\\u00E9s 0 %. Ajuste por CER si inflaci\\u00F3n > 10 %.\"],\n  [\"heading2\", \"4. Vencimiento\"],\n  [null, \"Si a 24 meses no hay ronda, el Prestamista puede exigir pago o conversi\\u00F3n a valor de mercado.\"],\n  [\"heading2\", \"5. Origen de Fondos\"],\n  [null, \"El Prestamista declara fondos l\\u00EDcitos (Ley 25.246).\"],\n  [\"heading2\", \"6. Ley y Arbitraje\"],\n  [null, \"Leyes argentinas. Arbitraje de derecho en CABA.\"],\n];\n\nconst body = context.document.body;\n\n// Insert every paragraph first (all still \"Normal\" at this point), THEN go\n// back and stamp the heading styles. Setting the style immediately after\n// each insertParagraph would make the *next* inserted paragraph inherit that\n// heading's pPr (Word carries the formatting of the paragraph mark it split\n// from), leaving stray <w:pStyle> on the body-text paragraphs.\nconst paragraphs = items.map(([, text]) => body.insertParagraph(text, \"End\"));\n\nitems.forEach(([style], i) => {\n  if (style) {\n    paragraphs[i].styleBuiltIn = Word.BuiltInStyleName[style];\n  }\n});\n\nawait context.sync();\n", "ps1": "# Build the \"CONTRATO DE MUTUO CONVERTIBLE\" template body.\n# Each entry pairs an (optional) built-in style name with the paragraph text;\n# $null keeps the paragraph on the default \"Normal\" style.\n$items = @(\n    @{ style = \"Heading 1\"; text = \"CONTRATO DE MUTUO CONVERTIBLE\" },\n    @{ style = $null; text = \"Entre ______________________ (el \u201cPrestamista\u201d) y la Compa\u00f1\u00eda en formaci\u00f3n (el \u201cMutuario\u201d), con fecha ___ de __________ de 2025.\" },\n    @{ style = \"Heading 2\"; text = \"1. Monto y Desembolso\" },\n    @{ style = $null; text = \"El Prestamista entrega USD __________ al Mutuario.\" },\n    @{ style = \"Heading 2\"; text = \"2. Conversi\u00f3n Autom\u00e1tica\" },\n    @{ style = $null; text = \"Al cerrar una ronda \u2265 USD ________, el mutuo se convierte en acciones preferidas con 20 % de descuento o valuation cap USD ________, lo que convenga al Prestamista.\" },\n    @{ style = \"Heading 2\"; text = \"3. Inter\u00e9s y Ajuste\" },\n    @{ style = $null; text = \"Inter\u00e9s 0 %. Ajuste por CER si inflaci\u00f3n > 10 %.\" },\n    @{ style = \"Heading 2\"; text = \"4. Vencimiento\" },\n    @{ style = $null; text = \"Si a 24 meses no hay ronda, el Prestamista puede exigir pago o conversi\u00f3n a valor de mercado.\" },\n    @{ style = \"Heading 2\"; text = \"5. Origen de Fondos\" },\n    @{ style = $null; text = \"El Prestamista declara fondos l\u00edcitos (Ley 25.246).\" },\n    @{ style = \"Heading 2\"; text = \"6. Ley y Arbitraje\" },\n    @{ style = $null; text = \"Leyes argentinas. Arbitraje de derecho en CABA.\" }\n)\n\n$d = $word.ActiveDocument\n\n# Append every paragraph first (all still \"Normal\" at this point), THEN go\n# back and stamp the heading styles. Setting the style immediately after\n# adding each paragraph would make the NEXT appended paragraph inherit that\n# heading's formatting (Word carries the formatting of the paragraph mark it\n# split from), leaving a stray style on the body-text paragraphs.\n$paras = @()\nforeach ($item in $items) {\n    $p = $d.Content.Paragraphs.Add()\n    $p.Range.Text = $item.text\n    $paras += $p\n}\n\nfor ($i = 0; $i -lt $items.Count; $i++) {\n    if ($items[$i].style) {\n        $paras[$i].Range.Style = $items[$i].style\n    }\n}\n"}
